{"js": "// Rewrite the \"console output\" submission-instructions bullet to ask for a\n// Word document (instead of a plain .txt file), and collapse the\n// \"Assignment reflection\" bullet (which used to be split around a\n// \"_GoBack\" bookmark) into a single run. The \"_GoBack\" bookmark itself\n// moves from the reflection bullet to the middle of the rewritten\n// console-output bullet.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet outputPara = null;\nlet reflectionPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Text file of your console output\") !== -1) {\n    outputPara = p;\n  } else if (p.text.indexOf(\"Assignment reflection in the comm\") !== -1) {\n    reflectionPara = p;\n  }\n}\n\nif (!outputPara || !reflectionPara) {\n  throw new Error(\"Could not locate target paragraphs\");\n}\n\n// The old \"_GoBack\" bookmark currently sits inside the reflection bullet;\n// drop it so the text there can be merged into a single clean run.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Bullet: \"Text file of your console output.txt ...\" -----------------\nconst oldOutputText =\n  \"Text file of your console output.txt (no images or Word documentation)\";\nconst firstPart =\n  \"Instead of attaching a text file with your console output, please attach a Word document that \";\nconst secondPart =\n  \"includes both your console output (using courier new font) and example of what each of your shapes look like running.\";\n\nconst outputMatches = outputPara.search(oldOutputText, { matchCase: true });\noutputMatches.load(\"items\");\nawait context.sync();\noutputMatches.items[0].insertText(\n  firstPart + secondPart,\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Re-insert the \"_GoBack\" bookmark at the boundary between the two halves\n// of the rewritten sentence.\nconst boundaryMatches = outputPara.search(firstPart, { matchCase: true });\nboundaryMatches.load(\"items\");\nawait context.sync();\nconst boundaryEnd = boundaryMatches.items[0].getRange(\n  Word.RangeLocation.after\n);\nboundaryEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Bullet: \"Assignment reflection in the comm\" + ... -------------------\nconst reflectionText =\n  \"Assignment reflection in the comment area the Good, Bad, and Ugly\";\nconst reflectionMatches = reflectionPara.search(reflectionText, {\n  matchCase: true,\n});\nreflectionMatches.load(\"items\");\nawait context.sync();\nreflectionMatches.items[0].insertText(\n  reflectionText,\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Rewrite the \"console output\" submission-instructions bullet to ask for a\n# Word document (instead of a plain .txt file), and collapse the\n# \"Assignment reflection\" bullet (which used to be split around a\n# \"_GoBack\" bookmark) into a single run. The \"_GoBack\" bookmark itself\n# moves from the reflection bullet to the middle of the rewritten\n# console-output bullet.\n\n$d = $word.ActiveDocument\n\n# Locate the two target paragraphs by their current text (index-independent).\n$outputPara = $null\n$reflectionPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"Text file of your console output*\") {\n        $outputPara = $p\n    } elseif ($t -like \"Assignment reflection in the comm*\") {\n        $reflectionPara = $p\n    }\n}\n\n# Drop the old \"_GoBack\" bookmark (it currently sits inside the reflection\n# bullet); this doesn't disturb character offsets.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Bullet: \"Text file of your console output.txt ...\" -----------------\n$firstPart = \"Instead of attaching a text file with your console output, please attach a Word document that \"\n$secondPart = \"includes both your console output (using courier new font) and example of what each of your shapes look like running.\"\n\n$outputRange = $outputPara.Range\n$outputStart = $outputRange.Start\n$outputEnd = $outputRange.End\n\n$replaceRange = $d.Range($outputStart, $outputEnd)\n$replaceRange.Text = $firstPart + $secondPart\n\n# Re-insert the \"_GoBack\" bookmark at the boundary between the two halves\n# of the rewritten sentence (a zero-length range right after $firstPart).\n$bookmarkPos = $outputStart + $firstPart.Length\n$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# --- Bullet: \"Assignment reflection in the comm\" + ... -------------------\n$reflectionText = \"Assignment reflection in the comment area the Good, Bad, and Ugly\"\n$reflectionRange = $reflectionPara.Range\n$reflectionReplace = $d.Range($reflectionRange.Start, $reflectionRange.End)\n$reflectionReplace.Text = $reflectionText\n"}
